$d = $word.ActiveDocument

function SafeReplace($paragraphRange, $old, $new) {
    $r1 = $paragraphRange.Duplicate
    $found = $r1.Find.Execute($old)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
        return $false
    }
    $oldStart = $r1.Start
    $oldEnd = $r1.End
    if ($new -ne "") {
        $r2 = $d.Range($oldEnd, $oldEnd)
        $r2.InsertBefore($new)
    }
    $r3 = $d.Range($oldStart, $oldEnd)
    $r3.Text = ""
    return $true
}


# Paragraph 1
$d.Hyperlinks.Item(1).TextToDisplay = '영어'
SafeReplace $d.Paragraphs(1).Range ' / Portuguese / French / Thai / Vietnamese / Spanish' ' / 포르투갈어 / 프랑스어 / 태국어 / 베트남어 / 스페인어' | Out-Null

# Paragraph 3
SafeReplace $d.Paragraphs(3).Range 'English' '영어' | Out-Null

# Paragraph 5
SafeReplace $d.Paragraphs(5).Range 'Brief' '개요' | Out-Null

# Paragraph 6
SafeReplace $d.Paragraphs(6).Range 'An email sent to partners in the target country who have sent their documents for review. It will be sent via customer.io' '검토를 위해 문서를 제출한 대상 국가의 파트너에게 보낼 이메일입니다. 고객.io를 통해 발송됩니다.' | Out-Null

# Paragraph 8
SafeReplace $d.Paragraphs(8).Range 'Target audience' '대상 청중' | Out-Null

# Paragraph 9
SafeReplace $d.Paragraphs(9).Range 'Invited partners who have submitted their documents' '문서를 제출한 초청 파트너' | Out-Null

# Paragraph 12
SafeReplace $d.Paragraphs(12).Range 'Subject line' '제목' | Out-Null
SafeReplace $d.Paragraphs(12).Range ' — we got your docs!  ' ' — 귀하의 문서를 받았습니다!  ' | Out-Null

# Paragraph 14
SafeReplace $d.Paragraphs(14).Range 'Thank you for submitting your documents' '귀하의 문서를 제출해 주셔서 감사합니다' | Out-Null

# Paragraph 16
SafeReplace $d.Paragraphs(16).Range ', ' ' ' | Out-Null
SafeReplace $d.Paragraphs(16).Range 'Hi ' '안녕하세요, ' | Out-Null

# Paragraph 18
SafeReplace $d.Paragraphs(18).Range 'Thank you for providing us with your documents for the upcoming ' '다가오는 ' | Out-Null
SafeReplace $d.Paragraphs(18).Range '. Based on the information you’ve given us, we’ll make the necessary arrangements, including accommodation and transportation.' '에 대한 귀하의 문서를 제공해 주셔서 감사합니다. 제공해 주신 정보를 바탕으로 숙박 및 교통편을 포함하여 필요한 사항을 준비해 드리겠습니다.' | Out-Null

# Paragraph 19
SafeReplace $d.Paragraphs(19).Range 'We’re currently reviewing your documents and will reach out to you if we need anything else. ' '저희는 현재 귀하의 문서를 검토 중이며 추가로 필요한 사항이 있으면 연락드리겠습니다. ' | Out-Null

# Paragraph 20
SafeReplace $d.Paragraphs(20).Range 'If you have any questions, please contact us via ' '궁금하신 점이 있는 경우, 저희 웹사이트의 ' | Out-Null
$d.Hyperlinks.Item(2).TextToDisplay = '실시간 채팅'
SafeReplace $d.Paragraphs(20).Range ' or ' ' 또는 ' | Out-Null
SafeReplace $d.Paragraphs(20).Range '. ' '을 통해 문의해 주시기 바랍니다. ' | Out-Null

# Paragraph 21
SafeReplace $d.Paragraphs(21).Range 'If you have any questions, please contact your country manager, ' '궁금하신 사항은, ' | Out-Null
SafeReplace $d.Paragraphs(21).Range ', at ' '에게 ' | Out-Null
SafeReplace $d.Paragraphs(21).Range ' or ' ' 또는 ' | Out-Null
SafeReplace $d.Paragraphs(21).Range ' (WhatsApp). ' ' (WhatsApp)을 통해 연락해 주시기 바랍니다. ' | Out-Null

# Paragraph 22
SafeReplace $d.Paragraphs(22).Range 'We look forward to seeing you at ' '' | Out-Null
SafeReplace $d.Paragraphs(22).Range '. ' '에서 만나 뵙기를 기대합니다. ' | Out-Null

# Comment text
SafeReplace $d.Comments(1).Range 'choose either one' '하나를 선택하세요' | Out-Null

